# Hortaliza, Comercializadora del Agro de Limarí - Poroto granado
# Weekly price update: insert a new data row for the latest week (2022-01-06)
# at row 36, pushing the existing rows 36-62 down to 37-63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 36; existing row 36..62 shift to 37..63
# (row formatting/number-format of the row above carries down as Excel does
# on a normal row insert).
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with this week's record.
$ws.Cells.Item(36, 1).Value2 = 2
$ws.Cells.Item(36, 2).Value2 = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(36, 3).Value2 = "Coquimbo"
$ws.Cells.Item(36, 4).Value2 = 44567
$ws.Cells.Item(36, 5).Value2 = 4
$ws.Cells.Item(36, 6).Value2 = 100112030
$ws.Cells.Item(36, 7).Value2 = "Poroto granado"
$ws.Cells.Item(36, 8).Value2 = "Sin especificar"
$ws.Cells.Item(36, 9).Value2 = "Primera"
$ws.Cells.Item(36, 10).Value2 = 360
$ws.Cells.Item(36, 11).Value2 = 15000
$ws.Cells.Item(36, 12).Value2 = 16000
$ws.Cells.Item(36, 13).Value2 = 15500
$ws.Cells.Item(36, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(36, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(36, 16).Value2 = 620
$ws.Cells.Item(36, 17).Value2 = 25
$ws.Cells.Item(36, 18).Value2 = "Hortaliza"
